$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("V4")

# New rows of 8-2 data (rows 102-106)
$ws.Range("A102").Value = "Rail 93650944 (after turn)"
$ws.Range("B102").Value = 30698
$ws.Range("C102").Value = 35992

$ws.Range("A103").Value = "Black screens"
$ws.Range("B103").Value = 30992
$ws.Range("C103").Value = 36290

$ws.Range("A104").Value = "Black screens"
$ws.Range("B104").Value = 31261
$ws.Range("C104").Value = 36576

$ws.Range("A105").Value = "Black screen"
$ws.Range("B105").Value = 31543
$ws.Range("C105").Value = 36872

$ws.Range("A106").Value = "Black screen (water scene)"
$ws.Range("B106").Value = 31909
$ws.Range("C106").Value = 37235

# Fill the Diff formula down for the 5 new rows in one shared operation
$ws.Range("D102:D106").Formula = "=IF(B102 >  0,C102-B102, 0)"

# New cells I93/J93
$ws.Range("I93").Value = 31543
$ws.Range("J93").Value = 36872

$ws.Range("J94").Select()


